# QRET_MechBOM_2019_20.xlsx - Add the upper avionics bulkhead eyebolt/nut rows
# and clarify the existing lower-bulkhead eyebolt & nut notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avionics")

# ---------------------------------------------------------------------------
# Row 20 - new upper-avionics-bulkhead eyebolt line item
#   (the string set first reclaims the shared-string slot orphaned by the
#    N19 note rewrite further below, so the table lands on the same index
#    order the real edit produced)
# ---------------------------------------------------------------------------
$ws.Range("C19").Copy($ws.Range("C20"))
$ws.Range("C20").Value = "5/16""-18 x 1-1/8"" thread length, Steel Eyebolt with shoulder"

$ws.Range("E20").Value = "McMaster-Carr"

$ws.Range("F19").Copy($ws.Range("F20"))
$ws.Range("F20").Value = "https://www.mcmaster.com/3014t46"

$ws.Range("G20").Value = "No"

$ws.Range("H19").Copy($ws.Range("H20"))
$ws.Range("H20").Value = 3.45

$ws.Range("I20").Value = 1

$ws.Range("N19").Copy($ws.Range("N20"))
$ws.Range("N20").Value = "This eyebolt is for the upper avionics bulkhead. It is rated for 900 lbs."

$ws.Rows.Item(20).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Row 21 - new upper-avionics-bulkhead nut line item
# ---------------------------------------------------------------------------
$ws.Range("C19").Copy($ws.Range("C21"))
$ws.Range("C21").Value = "Medium-Strength Steel Hex Nut - Grade 5, Zinc-Plated, 5/16""-18"

$ws.Range("E21").Value = "McMaster-Carr"

$ws.Range("F19").Copy($ws.Range("F21"))
$ws.Range("F21").Value = "https://www.mcmaster.com/95462a030"

$ws.Range("G21").Value = "No"

$ws.Range("H19").Copy($ws.Range("H21"))
$ws.Range("H21").Value = 7.18

$ws.Range("I21").Value = 100

$ws.Range("N19").Copy($ws.Range("N21"))

$ws.Rows.Item(21).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Row 19 - clarify the note: append the red "Check Home Depot..." callout
# (must happen AFTER the N19->N20/N21 format copies above, and before the
#  N21 text below, to keep the shared-string append order correct)
# ---------------------------------------------------------------------------
$ws.Range("N19").Value = "We need one nut to fasten the lower avionics bulkhead eyebolt **Check Home Depot - consider washer"
$note19 = $ws.Range("N19").Characters(65, 34)
$note19.Font.Color = 255
$note19.Font.Name = "Calibri"
$note19.Font.Size = 11

# ---------------------------------------------------------------------------
# Row 21 note text with red "Home Depot..." callout
# ---------------------------------------------------------------------------
$ws.Range("N21").Value = "Nut required for the upper avionics bulkhead. **Home Depot - consider washer"
$note21 = $ws.Range("N21").Characters(49, 28)
$note21.Font.Color = 255
$note21.Font.Name = "Calibri"
$note21.Font.Size = 11

# ---------------------------------------------------------------------------
# Hyperlinks for the two new supplier links (re-apply the normal hyperlink
# cell style afterwards, since Hyperlinks.Add stamps its own style index)
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.mcmaster.com/3014t46")
$ws.Range("F20").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F21"), "https://www.mcmaster.com/95462a030")
$ws.Range("F21").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# View state - scroll / selection moved down as more rows were added
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("N25").Select()
